# Update BunkerPrices at 2025-04-14 03:55
# Appends a new data row (row 32) to the bunker prices table, mirroring the
# values captured for the "Date" column's 2025-04-11 reading, and fixes up
# the "Date" column's number format so that only the newest row keeps the
# date-only format while the previous last row switches to the date-time
# format (matching the historical pattern used for earlier rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> value for the new row (row 32), in sheet column order.
$rowData = [ordered]@{
    "A"  = 565
    "B"  = 479
    "C"  = 450
    "D"  = 548
    "E"  = 498
    "F"  = 532
    "G"  = 474
    "H"  = 567
    "I"  = 495
    "J"  = 450
    "K"  = 572
    "L"  = 483
    "M"  = 462
    "N"  = 505
    "O"  = 557
    "P"  = 483
    "Q"  = 618
    "R"  = 497
    "S"  = 474
    "T"  = 485
    "U"  = 619
    "V"  = 540
    "W"  = 594
    "X"  = 490
    "Y"  = 45758
    "Z"  = 818
    "AA" = 556
    "AB" = 525.5
    "AC" = 512
    "AD" = 542
    "AE" = 504
    "AF" = 506
    "AG" = 730
    "AH" = 467
    "AI" = 720
    "AJ" = 474
    "AK" = 484
    "AL" = 550
    "AM" = 540
    "AN" = 485
    "AO" = 540
    "AP" = 529
    "AQ" = 565
    "AR" = 546
    "AS" = 635
    "AT" = 637
    "AU" = 493
    "AV" = 475
}

$newRow = 32

foreach ($col in $rowData.Keys) {
    $ws.Range("$col$newRow").Value = $rowData[$col]
}

# The new row's Date cell keeps the date-only format (matches style used by
# every other row before the previous last row).
$ws.Range("Y$newRow").NumberFormat = "YYYY-MM-DD"

# The row that used to be last (31) switches to the date-time format, as
# happens for the row preceding the newest entry.
$ws.Range("Y31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
